$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L2").Value = "[60.32903324416962, 67.09456835999733]"
$ws.Range("T2").Value = "[47.417751145959606, 52.129590525942525]"
$ws.Range("L3").Value = "[58.07780709576167, 68.71255590107361]"
$ws.Range("T3").Value = "[47.12377482695696, 53.28545000888742]"
